$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "SNO"
$ws.Range("B1").Value = "OUR CODE"
$ws.Range("C1").Value = "OUT"
$ws.Range("D1").Value = "OUT"
$ws.Range("E1").Value = "OUT"
$ws.Range("F1").Value = "IN"

# --- Row 2 (blank SNO/code, text-like dates) ---
$ws.Range("A2:F2").NumberFormat = "@"
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "30/10/2024"
$ws.Range("D2").Value = "30/9/2024"
$ws.Range("E2").Value = "30/8/2024"
$ws.Range("F2").Value = "1/8/24"

# --- Rows 3-5 (numbers stored as text) ---
$ws.Range("A3:F5").NumberFormat = "@"

$ws.Range("A3").Value = "1"
$ws.Range("B3").Value = "901"
$ws.Range("C3").Value = "21"
$ws.Range("D3").Value = "37"
$ws.Range("E3").Value = "64"
$ws.Range("F3").Value = "295"

$ws.Range("A4").Value = "2"
$ws.Range("B4").Value = "902"
$ws.Range("C4").Value = "28"
$ws.Range("D4").Value = "25"
$ws.Range("E4").Value = "46"
$ws.Range("F4").Value = "268"

$ws.Range("A5").Value = "3"
$ws.Range("B5").Value = "903"
$ws.Range("C5").Value = "16"
$ws.Range("D5").Value = "45"
$ws.Range("E5").Value = "41"
$ws.Range("F5").Value = "248"

# --- Remove the old G/H columns (DATE / NOTES) ---
$ws.Range("H1:H3").EntireColumn.Delete()
$ws.Range("G1:G3").EntireColumn.Delete()

# --- Column widths (Excel stores width + ~0.8333 padding, so compensate) ---
$pad = 5/6
$ws.Columns.Item(1).ColumnWidth = 8 - $pad
$ws.Columns.Item(2).ColumnWidth = 15 - $pad
$ws.Columns.Item(3).ColumnWidth = 12 - $pad
$ws.Columns.Item(4).ColumnWidth = 12 - $pad
$ws.Columns.Item(5).ColumnWidth = 12 - $pad
$ws.Columns.Item(6).ColumnWidth = 12 - $pad
